$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph that currently reads:
#   "Complete Project is very big to be uploaded here. Therefore, I uploaded
#    it in github you can access it via this link"
# and the paragraph right after it (which only holds the _GoBack bookmark).
# We look the text up dynamically instead of hard-coding paragraph indexes so
# the script keeps working even if earlier paragraphs shift slightly.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$noteIdx = 0
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*uploaded it in*you can access it via this link*") {
        $noteIdx = $i
    }
}

$notePara = $d.Paragraphs($noteIdx).Range
$bookmarkPara = $d.Paragraphs($noteIdx + 1).Range

$rangeStart = $notePara.Start
$rangeEnd = $bookmarkPara.End

$driveUrl = "https://drive.google.com/drive/folders/1K2y6_13NYsaSQ8uiFGceDRon151fOYqi?usp=sharing"

# ---------------------------------------------------------------------------
# Replace the note paragraph + the bookmark-only paragraph with:
#   1. the same note paragraph, but "github" swapped for "Google drive" (and
#      the spell-check proofErr markers that wrapped it gone), followed by
#      the _GoBack bookmark moved to the end of that paragraph;
#   2. a brand new paragraph containing the Google Drive link as a hyperlink;
#   3. a trailing empty paragraph (matching formatting) before the sectPr.
# Using InsertXML gives precise control of run/paragraph structure, which is
# needed to drop the <w:proofErr/> marks and to relocate the bookmark.
# ---------------------------------------------------------------------------
$r = $d.Range($rangeStart, $rangeEnd)
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
<w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve">Complete Project is very big to be uploaded here. Therefore, I uploaded it in </w:t></w:r>
<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr><w:t>Google drive</w:t></w:r>
<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr><w:t xml:space="preserve"> you can access it via this link</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
<w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr></w:pPr>
<w:hyperlink r:id="rIdGoogleDrivePlaceholder" w:tooltip="$driveUrl" w:history="1">
<w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr><w:t>$driveUrl</w:t></w:r>
</w:hyperlink>
</w:p>
<w:p>
<w:pPr><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-IN"/></w:rPr></w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rIdGoogleDrivePlaceholder" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="$driveUrl" TargetMode="External"/>
</Relationships>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$r.InsertXML($xml)

# ---------------------------------------------------------------------------
# InsertXML does not keep character-style references (w:rStyle), so the new
# hyperlink run lost its "Hyperlink" style during the insert above. Restore
# it the same way real Word would apply the style when you create a
# hyperlink: set the Style on the run's own range (the run is the only thing
# in its paragraph, besides the paragraph mark).
# ---------------------------------------------------------------------------
$newHyperlinkParaIdx = $noteIdx + 1
$hyperlinkRunRange = $d.Paragraphs($newHyperlinkParaIdx).Range
$hyperlinkRunRange.MoveEnd(1, -1) | Out-Null
$hyperlinkRunRange.Style = "Hyperlink"

Write-Host "Google Drive link added."
